# ERP-521 - As a Case Work Clerk I want the outstation address to appear
# on letters generated for Scottish outstations (Aberdeen, Dundee, Edinburgh)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New data rows, in the exact order they must be written so that the
# shared-string table is built up in the same sequence as the target
# workbook (first-use order == index order).
# Columns: Field name (A) / Value (B) / wrap style flag (copy from B14
# when true, copy from B16 -- plain style -- when false) / hyperlink
# mailto address (only set on the "...Email" rows).
# ---------------------------------------------------------------------
$rows = @(
    @{ Row=21; Field="tribunalAberdeenAddressLine1"; Value="Ground Floor";               Wrap=$true },
    @{ Row=22; Field="tribunalAberdeenAddressLine2"; Value="AB1, 48 Huntly Street";       Wrap=$true },
    @{ Row=23; Field="tribunalAberdeenTown";         Value="Aberdeen";                    Wrap=$true },
    @{ Row=24; Field="tribunalAberdeenPostCode";     Value="AB10 1SH";                    Wrap=$true },
    @{ Row=25; Field="tribunalAberdeenTelephone";    Value="01224 593 137";               Wrap=$false },
    @{ Row=26; Field="tribunalAberdeenFax";          Value="0870 761 7766";               Wrap=$false },
    @{ Row=27; Field="tribunalAberdeenDX";           Value="DX AB77";                     Wrap=$false },
    @{ Row=28; Field="tribunalAberdeenEmail";        Value="aberdeenet@justice.gov.uk";    Wrap=$false; Mail=$true },

    @{ Row=29; Field="tribunalDundeeAddressLine1";   Value="Ground Floor";                Wrap=$true },
    @{ Row=30; Field="tribunalDundeeAddressLine2";   Value="Block C, Caledonian House";   Wrap=$true },
    @{ Row=31; Field="tribunalDundeeAddressLine3";   Value="Greenmarket";                 Wrap=$true },
    @{ Row=32; Field="tribunalDundeeTown";           Value="Dundee";                      Wrap=$true },
    @{ Row=33; Field="tribunalDundeePostCode";       Value="DD1 4QG";                     Wrap=$true },
    @{ Row=34; Field="tribunalDundeeTelephone";      Value="01382 221 578";               Wrap=$false },
    @{ Row=35; Field="tribunalDundeeFax";            Value="01382 227 136";               Wrap=$false },
    @{ Row=36; Field="tribunalDundeeDX";             Value="DX DD51";                     Wrap=$false },
    @{ Row=37; Field="tribunalDundeeEmail";          Value="dundeeet@justice.gov.uk";      Wrap=$false; Mail=$true },

    @{ Row=38; Field="tribunalEdinburghAddressLine1"; Value="54-56 Melville Street";      Wrap=$true },
    @{ Row=39; Field="tribunalEdinburghTown";         Value="Edinburgh";                  Wrap=$true },
    @{ Row=40; Field="tribunalEdinburghPostCode";     Value="EH3 7HF";                    Wrap=$true },
    @{ Row=41; Field="tribunalEdinburghTelephone";    Value="0131 226 5584";              Wrap=$false },
    @{ Row=42; Field="tribunalEdinburghFax";          Value="0131 220 6847";              Wrap=$false },
    @{ Row=43; Field="tribunalEdinburghDX";           Value="DX ED147";                   Wrap=$false },
    @{ Row=44; Field="tribunalEdinburghEmail";        Value="edinburghet@justice.gov.uk";  Wrap=$false; Mail=$true }
)

# Reference cells already carrying the two formats used throughout the
# sheet -- copy their formatting instead of inventing new style entries.
$wrapFormatSource  = $ws.Range("B14")   # s="4" (wrap / address-style) cells
$plainFormatSource = $ws.Range("B16")   # s="0" (plain) cells

foreach ($r in $rows) {
    $ws.Range("A" + $r.Row).Value = $r.Field
    $ws.Range("B" + $r.Row).Value = $r.Value

    if ($r.Wrap) {
        $wrapFormatSource.Copy()
    } else {
        $plainFormatSource.Copy()
    }
    $ws.Range("B" + $r.Row).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# Hyperlinks for the three new email rows
$ws.Hyperlinks.Add($ws.Range("B28"), "mailto:aberdeenet@justice.gov.uk", [System.Type]::Missing, [System.Type]::Missing, "aberdeenet@justice.gov.uk")
$ws.Range("B14").Copy()
$ws.Range("B28").PasteSpecial(-4122)
$plainFormatSource.Copy()
$ws.Range("B28").PasteSpecial(-4122)

$ws.Hyperlinks.Add($ws.Range("B37"), "mailto:dundeeet@justice.gov.uk", [System.Type]::Missing, [System.Type]::Missing, "dundeeet@justice.gov.uk")
$plainFormatSource.Copy()
$ws.Range("B37").PasteSpecial(-4122)

$ws.Hyperlinks.Add($ws.Range("B44"), "mailto:edinburghet@justice.gov.uk", [System.Type]::Missing, [System.Type]::Missing, "edinburghet@justice.gov.uk")
$plainFormatSource.Copy()
$ws.Range("B44").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Move the viewport / selection down to the newly added Scottish block,
# matching the author's saved cursor position.
$ws.Range("A20:B44").Select()
$excel.ActiveWindow.ScrollRow = 20
$excel.ActiveWindow.ScrollColumn = 1
